$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '27.684.48'
Set-TextValue "E2" '  +0.01%  '
Set-TextValue "D3" '1.634.49'
Set-TextValue "E3" '  -0.26%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '212.16'
Set-TextValue "E5" '  -0.14%  '
Set-TextValue "D6" '0.523'
Set-TextValue "E6" '  -0.31%  '
Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  -0.04%  '
Set-TextValue "D8" '23.35'
Set-TextValue "E8" '  +1.06%  '
Set-TextValue "E9" '  +2.26%  '
Set-TextValue "E10" '  +0.16%  '
Set-TextValue "D11" '0.0860'
Set-TextValue "E11" '  -3.74%  '
Set-TextValue "D12" '1.866.24'
Set-TextValue "E12" '  -0.26%  '
Set-TextValue "D13" '1.641.28'
Set-TextValue "E13" '  +0.16%  '
Set-TextValue "E14" '  -0.32%  '
Set-TextValue "E15" '  -1.24%  '
Set-TextValue "D16" '65.12'
Set-TextValue "D17" '27.664.49'
Set-TextValue "E17" '  -0.01%  '
Set-TextValue "D18" '230.53'
Set-TextValue "E18" '  -0.23%  '
Set-TextValue "E19" '  -0.29%  '
Set-TextValue "E20" '  -1.54%  '
Set-TextValue "D21" '0.999'
Set-TextValue "E21" '  -0.04%  '
Set-TextValue "D22" '10.60'
Set-TextValue "E22" '  +3.81%  '
Set-TextValue "E23" '  +1.01%  '
Set-TextValue "E24" '  +3.33%  '
Set-TextValue "D25" '148.73'
Set-TextValue "E25" '  -1.50%  '
Set-TextValue "D26" '6.89'
Set-TextValue "E26" '  -1.04%  '
Set-TextValue "E27" '  -0.31%  '
Set-TextValue "B28" 'EthereumClassic'
Set-TextValue "C28" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D28" '15.54'
Set-TextValue "E28" '  -0.44%  '
Set-TextValue "B29" 'BinanceUSD'
Set-TextValue "C29" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D29" '1.00'
Set-TextValue "E29" '  +0.00%  '
Set-TextValue "E30" '  +0.17%  '
Set-TextValue "E32" '  -1.06%  '
Set-TextValue "D33" '1.481.93'
Set-TextValue "E33" '  +1.63%  '
Set-TextValue "D34" '3.09'
Set-TextValue "E34" '  -1.16%  '
Set-TextValue "E35" '  -2.04%  '
Set-TextValue "E36" '  -1.41%  '
Set-TextValue "D37" '0.960'
Set-TextValue "E37" '  +6.94%  '
Set-TextValue "E38" '  +0.32%  '
Set-TextValue "D39" '0.558'
Set-TextValue "E39" '  -1.39%  '
Set-TextValue "E40" '  +0.08%  '
Set-TextValue "E41" '  +0.96%  '
Set-TextValue "E42" '  -0.03%  '
Set-TextValue "D43" '67.67'
Set-TextValue "E43" '  -2.99%  '
Set-TextValue "D44" '2.46'
Set-TextValue "E44" '  -0.22%  '
Set-TextValue "D45" '2.21'
Set-TextValue "E45" '  -1.32%  '
Set-TextValue "E46" '  -5.10%  '
Set-TextValue "D47" '1.775.46'
Set-TextValue "E47" '  -0.33%  '
Set-TextValue "E48" '  +0.63%  '
Set-TextValue "D49" '87.58'
Set-TextValue "E49" '  +1.03%  '
Set-TextValue "E50" '  -1.73%  '
Set-TextValue "E51" '  -0.12%  '
